$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("F1").Value = "Contact No."
$ws.Range("G1").Value = "Fee status"

# Row 2 - replace Chetan Kaushik's record with Vaibhav Aggarwal's
$ws.Range("B2").Value = "Vaibhav Aggarwal"
$ws.Range("C2").Value = "Computer Engineering"
$ws.Range("D2").Value = "vaibhav.aggarwal.@computer.jcboseust.ac.in"
$ws.Range("E2").Value = "&negZrBrA8?"
$ws.Range("F2").Value = 9764767579
$ws.Range("G2").Value = 45000

# Row 3 - replace Dhruv Singh's record with Vidushi Tickoo's
$ws.Range("B3").Value = "Vidushi Tickoo"
$ws.Range("C3").Value = "Computer Engineering"
$ws.Range("D3").Value = "vidushi.tickoo.@computer.jcboseust.ac.in"
$ws.Range("E3").Value = "w:,^ROykm!|"
$ws.Range("F3").Value = 8920021900
$ws.Range("G3").Value = 45000
